# Final commit of upload excel file:
# - update a few contact field values (first names, street, hobby lists)
# - bump the header/data row heights slightly (18.75 -> 19.5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (tintu / S / ... ) - FirstName tintu -> rohan, Hobbies trailing comma removed
$ws.Range("B2").Value = "rohan"
$ws.Range("L2").Value = "Reading ,Drawing"

# Row 3 (Maya / S / ... ) - FirstName Maya -> mini, Street dfbdf -> abcd,
# Hobbies trailing comma removed
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# Rows 1-3 got a touch taller in the re-exported sheet.
$ws.Range("1:3").RowHeight = 19.5
